# "merge additional data into original dataset"
# Normalizes the free-text reason labels used across reviewers and adds a
# consolidated "Decision" column (E) that captures the final, merged call
# for each paper.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("exclusion results")

# --- Normalize inconsistent reviewer wording to the canonical labels -------
$ws.Range("D2").Value  = "Philosophical, Opinion, or Experience paper"
$ws.Range("C3").Value  = "No MDE application"
$ws.Range("C4").Value  = "No MDE application"
$ws.Range("C5").Value  = "No MDE application"
$ws.Range("C6").Value  = "No MDE application"
$ws.Range("C9").Value  = "No MDE application"
$ws.Range("C10").Value = "No MDE application"
$ws.Range("C11").Value = "No MDE application"
$ws.Range("C12").Value = "No MDE application"
$ws.Range("C13").Value = "No MDE application"
$ws.Range("B14").Value = "Philosophical, Opinion, or Experience paper"
$ws.Range("C14").Value = "Philosophical, Opinion, or Experience paper"
$ws.Range("D14").Value = "Philosophical, Opinion, or Experience paper"
$ws.Range("B18").Value = "Philosophical, Opinion, or Experience paper"
$ws.Range("C18").Value = "No MDE application"

# --- Add the merged "Decision" column --------------------------------------
$ws.Range("E1").Value  = "Decision"
$ws.Range("E2").Value  = "No MDE application"
$ws.Range("E3").Value  = "No MDE application"
$ws.Range("E4").Value  = "No MDE application"
$ws.Range("E5").Value  = "No MDE application"
$ws.Range("E6").Value  = "No MDE application"
$ws.Range("E7").Value  = "No DT used"
$ws.Range("E9").Value  = "No MDE application"
$ws.Range("E10").Value = "No MDE application"
$ws.Range("E11").Value = "No MDE application"
$ws.Range("E12").Value = "No MDE application"
$ws.Range("E13").Value = "No MDE application"
$ws.Range("E14").Value = "Philosophical, Opinion, or Experience paper"
$ws.Range("E15").Value = "No MDE application"
$ws.Range("E16").Value = "No DT used"
$ws.Range("E17").Value = "No DT used"
$ws.Range("E18").Value = "No MDE application"
$ws.Range("E20").Value = "No MDE application"
$ws.Range("E24").Value = "No MDE application"

# --- Reflect the reviewer's last on-screen position/selection --------------
$ws.Range("D2:D24").Select()
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Application.ActiveWindow.ScrollColumn = 3
